$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "26.068.23"
Set-TextValue $ws.Range("E2") "  -0.88%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.649.65"
Set-TextValue $ws.Range("E3") "  -0.88%  "

# Row 5
Set-TextValue $ws.Range("D5") "217.44"
Set-TextValue $ws.Range("E5") "  -0.73%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.5210"
Set-TextValue $ws.Range("E6") "  -2.37%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.35%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.2615"
Set-TextValue $ws.Range("E8") "  -1.56%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.06282"
Set-TextValue $ws.Range("E9") "  -1.81%  "

# Row 10
Set-TextValue $ws.Range("D10") "20.48"
Set-TextValue $ws.Range("E10") "  -0.53%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.07793"
Set-TextValue $ws.Range("E11") "  -0.43%  "

# Row 12
Set-TextValue $ws.Range("D12") "4.476"
Set-TextValue $ws.Range("E12") "  -2.05%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.648.75"
Set-TextValue $ws.Range("E13") "  -1.16%  "

# Row 14
Set-TextValue $ws.Range("D14") "1.877.63"
Set-TextValue $ws.Range("E14") "  -0.77%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.5531"
Set-TextValue $ws.Range("E15") "  +0.15%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.0₅7988"
Set-TextValue $ws.Range("E16") "  -2.77%  "

# Row 17
Set-TextValue $ws.Range("D17") "64.72"
Set-TextValue $ws.Range("E17") "  -1.54%  "

# Row 18
Set-TextValue $ws.Range("D18") "26.061.56"
Set-TextValue $ws.Range("E18") "  -0.94%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -1.41%  "

# Row 21
Set-TextValue $ws.Range("D21") "193.94"
Set-TextValue $ws.Range("E21") "  -0.18%  "

# Row 22
Set-TextValue $ws.Range("E22") "  -1.58%  "

# Row 23
Set-TextValue $ws.Range("D23") "5.940"
Set-TextValue $ws.Range("E23") "  -1.67%  "

# Row 24
Set-TextValue $ws.Range("E24") "  -0.37%  "

# Row 25
Set-TextValue $ws.Range("D25") "146.65"
Set-TextValue $ws.Range("E25") "  +0.34%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.1201"
Set-TextValue $ws.Range("E26") "  -2.76%  "

# Row 27
Set-TextValue $ws.Range("D27") "7.171"
Set-TextValue $ws.Range("E27") "  -0.23%  "

# Row 28
Set-TextValue $ws.Range("E28") "  -1.39%  "

# Row 29
Set-TextValue $ws.Range("E29") "  -0.45%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.05596"
Set-TextValue $ws.Range("E30") "  -4.31%  "

# Row 31
Set-TextValue $ws.Range("E31") "  -1.08%  "

# Row 32
Set-TextValue $ws.Range("D32") "3.480"
Set-TextValue $ws.Range("E32") "  -3.82%  "

# Row 33
Set-TextValue $ws.Range("D33") "3.357"
Set-TextValue $ws.Range("E33") "  +2.26%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.587"
Set-TextValue $ws.Range("E34") "  -1.58%  "

# Row 35
Set-TextValue $ws.Range("E35") "  -1.09%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.9473"
Set-TextValue $ws.Range("E36") "  -1.68%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.413"
Set-TextValue $ws.Range("E37") "  -0.16%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.5629"
Set-TextValue $ws.Range("E38") "  -3.00%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.01584"

# Row 40
Set-TextValue $ws.Range("E40") "  +1.13%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.057.60"
Set-TextValue $ws.Range("E41") "  +0.47%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -0.39%  "

# Row 43
Set-TextValue $ws.Range("E43") "  -3.17%  "

# Row 44
Set-TextValue $ws.Range("D44") "102.29"
Set-TextValue $ws.Range("E44") "  -2.16%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.789.07"
Set-TextValue $ws.Range("E45") "  -0.75%  "

# Row 46
Set-TextValue $ws.Range("D46") "57.08"
Set-TextValue $ws.Range("E46") "  -1.21%  "

# Row 47
Set-TextValue $ws.Range("E47") "  -0.46%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.05309"
Set-TextValue $ws.Range("E48") "  +2.81%  "

# Row 49
Set-TextValue $ws.Range("B49") "Mantle"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D49") "0.4332"
Set-TextValue $ws.Range("E49") "  -1.20%  "

# Row 50
Set-TextValue $ws.Range("B50") "EnergySwap"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.938"
Set-TextValue $ws.Range("E50") "  -1.23%  "

# Row 51
Set-TextValue $ws.Range("B51") "BabyDogeCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D51") "0.0₈101"
Set-TextValue $ws.Range("E51") "  -4.49%  "
